$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J ("I0" and "IF"), matching the
# formatting already used by the existing header row (copy format from H1,
# the last existing header cell, so the new headers reuse the same style
# rather than minting new ones).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF data values for rows 2-23
$data = @{
    2  = @(8, 9)
    3  = @(4, 4)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(6, 6)
    9  = @(8, 8)
    10 = @(6, 6)
    11 = @(8, 9)
    12 = @(7, 7)
    13 = @(3, 3)
    14 = @(10, 10)
    15 = @(3, 4)
    16 = @(8, 8)
    17 = @(7, 7)
    18 = @(9, 9)
    19 = @(7, 7)
    20 = @(8, 8)
    21 = @(6, 6)
    22 = @(6, 6)
    23 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
